# POP1.docx template adjustments:
#  - move the "_GoBack" bookmark so it sits inside "Natureza" (between
#    "Na" and "tureza...") instead of right before "monitoramento"
#  - split the two "(nome completo)" runs into three runs each, with
#    proofErr gramStart/gramEnd markers bracketing the word "nome"
#  - add left-indent to the Ttulo1 / Ttulo2 (heading 1 / heading 2) styles

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Relocate the "_GoBack" bookmark into the "Natureza..." heading,
#    right after the first two letters ("Na" | bookmark | "tureza...").
#    Re-adding a bookmark with the same name moves it (names are unique),
#    so the stale one near "monitoramento" disappears automatically.
# ---------------------------------------------------------------------
$findRange = $d.Content
$found = $findRange.Find.Execute(
    "Natureza da superfície a ser higienizada",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $splitPos = $findRange.Start + 2
    $bookmarkRange = $d.Range($splitPos, $splitPos)
    $d.Bookmarks.Add("_GoBack", $bookmarkRange)
}

# ---------------------------------------------------------------------
# 2) Split every "(nome completo)" run into "(" + "nome" + " completo)",
#    wrapping the middle run with gramStart/gramEnd proofErr markers -
#    matches Word re-flagging the word after a retype/grammar pass.
# ---------------------------------------------------------------------
$nomeCompletoXml = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:rFonts w:cs="Arial"/></w:rPr><w:t xml:space="preserve">(</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:cs="Arial"/></w:rPr><w:t>nome</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:cs="Arial"/></w:rPr><w:t xml:space="preserve"> completo)</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$scan = $d.Content
$scan.Start = 0
$scan.End = $d.Content.End
$guard = 0
while ($scan.Find.Execute("(nome completo)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    $target = $d.Range($scan.Start, $scan.End)
    $target.InsertXML($nomeCompletoXml)

    $scan.Start = $target.End
    $scan.End = $d.Content.End

    $guard = $guard + 1
    if ($guard -gt 10) { break }
}

# ---------------------------------------------------------------------
# 3) Give the Ttulo1 / Ttulo2 heading styles a left indent (426 / 567
#    twips == 21.3pt / 28.35pt) so the numbered headings line up with
#    the new margins.
# ---------------------------------------------------------------------
$ttulo1 = $d.Styles("Ttulo1")
$ttulo1.ParagraphFormat.LeftIndent = 21.3

$ttulo2 = $d.Styles("Ttulo2")
$ttulo2.ParagraphFormat.LeftIndent = 28.35
